# Apply the "previously unexpressing" update from Claire:
#  - two rows that had a label but no measurement data ("c17" and "n33")
#    are removed entirely (their constructs were never actually measured),
#    which shifts every row below them up.
#  - the surviving rows are relabeled: the "chimera" header becomes "name",
#    and the first two data rows (previously "CsChrim"/"C1C2") become
#    "cschrimson"/"c1c2".
#  - the active selection is left on C15 to match where Claire was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that only ever had a label in column A and no data
# (row 6 = "c17", row 24 = "n33"). Deleting row 6 first shifts the second
# one up from 24 to 23.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(23).Delete()

# Relabel the header and the first two (renamed) constructs.
$ws.Range("A1").Value = "name"
$ws.Range("A2").Value = "cschrimson"
$ws.Range("A3").Value = "c1c2"

# Leave the selection where Claire left it.
$ws.Range("C15").Select()
